# Renamed few transcripts. Updated the DataSheet
# Column D ("Speaker") contains tags that need to be shortened:
#   "RT1"      -> "T"
#   "Students" -> "SS"
#   "Class"    -> "SS"   (only the single occurrence at D25)
# Rows where D holds an actual person's name (e.g. "David", "Andrew",
# "Michael", "Brian", "Erik") must stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 101 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value2

    if ($val -eq "RT1") {
        $cell.Value = "T"
    }
    elseif ($val -eq "Students") {
        $cell.Value = "SS"
    }
    elseif ($val -eq "Class") {
        $cell.Value = "SS"
    }
}
